$d = $word.ActiveDocument

# Locate the "Planejamento do Projeto / Marcos do Projeto*" heading
# paragraph, then step back to the blank paragraph that immediately
# precedes it (two sz=28/szCs=22 runs followed by an empty run) -
# that blank paragraph is where the new divider line gets appended.
$headingIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("Planejamento do Projeto")) {
        $headingIndex = $i
        break
    }
}

if ($headingIndex -eq -1) {
    throw "heading paragraph not found"
}

$target = $d.Paragraphs.Item($headingIndex - 1)

if ($target -eq $null) {
    throw "target paragraph not found"
}

$r = $target.Range

# Rebuild the target paragraph exactly as-is, then append a brand-new
# paragraph (section divider made of em dashes) right after it. Using
# InsertXML on the paragraph's own range (rather than a collapsed point)
# keeps the new run free of any inherited rPr, matching a freshly
# authored divider line.
$xml = @"
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:pPr><w:jc w:val='both'/><w:spacing w:lineRule='auto' w:line='360'/></w:pPr><w:r><w:rPr><w:bCs w:val='0'/><w:sz w:val='28'/><w:szCs w:val='22'/></w:rPr></w:r><w:r><w:rPr><w:bCs w:val='0'/><w:sz w:val='28'/><w:szCs w:val='22'/></w:rPr></w:r><w:r/></w:p><w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:pPr><w:jc w:val='both'/><w:spacing w:lineRule='auto' w:line='360'/></w:pPr><w:r><w:t xml:space='preserve'>———————————————————————————————————</w:t></w:r><w:r/></w:p>
"@

[void]$r.InsertXML($xml)
